$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B2").Value = 0.00000000000000018732047832379676

$ws.Range("B3").Value = 8.6741028045005617
$ws.Range("C3").Value = 19.749328915738147
$ws.Range("D3").Value = 340.25067108426191

$ws.Range("B4").Value = 15.460914713021179
$ws.Range("C4").Value = 21.361084002807804
$ws.Range("D4").Value = 338.63891599719221

$ws.Range("B5").Value = 20.976379169051555
$ws.Range("C5").Value = 22.864417438491312
$ws.Range("D5").Value = 337.13558256150867

$ws.Range("B6").Value = 25.587849945264185
$ws.Range("C6").Value = 24.279515549106414
$ws.Range("D6").Value = 335.72048445089359

$ws.Range("B7").Value = 29.529777375310879
$ws.Range("C7").Value = 25.621039470848146
$ws.Range("D7").Value = 334.37896052915187

$ws.Range("B8").Value = 32.959746344041726
$ws.Range("C8").Value = 26.900034517825578
$ws.Range("D8").Value = 333.09996548217441

$ws.Range("B9").Value = 35.988105245180286
$ws.Range("C9").Value = 28.125066913004545
$ws.Range("D9").Value = 331.87493308699544

$ws.Range("B10").Value = 38.694736190202732
$ws.Range("C10").Value = 29.302938728976116
$ws.Range("D10").Value = 330.69706127102393

$ws.Range("B11").Value = 41.139080203554784
$ws.Range("C11").Value = 30.439157851663548
$ws.Range("D11").Value = 329.56084214833646
